$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --------------------------------------------------------------------------
# Refresh the "Price" (D) and "Volume(1h)" (E) columns with the latest
# coinranking.com snapshot values (scheduled GitHub Actions data refresh).
#
# D-column values are written with a leading apostrophe to force plain text
# (many values look numeric, e.g. "1.00", "6.20", "0.0670", and Excel would
# otherwise silently reinterpret/round them as numbers, dropping trailing
# zeros or the "thousands-as-dots" notation used for the larger coins).
# --------------------------------------------------------------------------

$ws.Range("D2").Value = "'57.944.17"
$ws.Range("E2").Value = "  +0.40%  "

$ws.Range("D3").Value = "'3.107.79"
$ws.Range("E3").Value = "  +1.69%  "

$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").Value = "'526.11"
$ws.Range("E5").Value = "  +2.50%  "

$ws.Range("D6").Value = "'142.58"
$ws.Range("E6").Value = "  +2.03%  "

$ws.Range("D8").Value = "'3.106.21"
$ws.Range("E8").Value = "  +1.65%  "

$ws.Range("E9").Value = "  +1.96%  "

$ws.Range("D10").Value = "'7.24"
$ws.Range("E10").Value = "  -1.39%  "

$ws.Range("E11").Value = "  +1.75%  "

$ws.Range("D12").Value = "'0.392"
$ws.Range("E12").Value = "  +4.58%  "

$ws.Range("D13").Value = "'3.637.65"
$ws.Range("E13").Value = "  +1.60%  "

$ws.Range("E14").Value = "  +1.96%  "

$ws.Range("D15").Value = "'25.68"
$ws.Range("E15").Value = "  -1.93%  "

$ws.Range("D16").Value = "'0.0000165"
$ws.Range("E16").Value = "  +2.01%  "

$ws.Range("D17").Value = "'58.001.25"
$ws.Range("E17").Value = "  +0.46%  "

$ws.Range("D18").Value = "'3.126.79"
$ws.Range("E18").Value = "  +2.05%  "

$ws.Range("D19").Value = "'6.11"
$ws.Range("E19").Value = "  -0.42%  "

$ws.Range("D20").Value = "'12.83"
$ws.Range("E20").Value = "  +0.73%  "

$ws.Range("E21").Value = "  -0.16%  "

$ws.Range("D22").Value = "'342.47"
$ws.Range("E22").Value = "  +3.82%  "

$ws.Range("E23").Value = "  +0.02%  "

$ws.Range("D24").Value = "'0.514"
$ws.Range("E24").Value = "  +3.35%  "

$ws.Range("D25").Value = "'67.32"
$ws.Range("E25").Value = "  +3.94%  "

$ws.Range("D26").Value = "'0.169"
$ws.Range("E26").Value = "  -1.29%  "

$ws.Range("D27").Value = "'1.01"
$ws.Range("E27").Value = "  +0.58%  "

$ws.Range("D28").Value = "'0.0₃0923"
$ws.Range("E28").Value = "  +3.15%  "

$ws.Range("D29").Value = "'6.48"
$ws.Range("E29").Value = "  +1.60%  "

$ws.Range("D30").Value = "'0.998"
$ws.Range("E30").Value = "  -0.06%  "

$ws.Range("D31").Value = "'7.29"
$ws.Range("E31").Value = "  +2.57%  "

$ws.Range("D32").Value = "'1.88"
$ws.Range("E32").Value = "  +4.61%  "

$ws.Range("D33").Value = "'21.08"
$ws.Range("E33").Value = "  +2.03%  "

$ws.Range("D34").Value = "'1.21"
$ws.Range("E34").Value = "  +1.96%  "

$ws.Range("D35").Value = "'158.11"
$ws.Range("E35").Value = "  +1.85%  "

$ws.Range("D36").Value = "'4.69"
$ws.Range("E36").Value = "  +3.79%  "

$ws.Range("D37").Value = "'6.20"
$ws.Range("E37").Value = "  +4.01%  "

$ws.Range("D38").Value = "'26.58"
$ws.Range("E38").Value = "  -2.25%  "

$ws.Range("D39").Value = "'1.26"
$ws.Range("E39").Value = "  +0.22%  "

$ws.Range("D40").Value = "'0.0670"
$ws.Range("E40").Value = "  -0.45%  "

$ws.Range("D41").Value = "'4.05"
$ws.Range("E41").Value = "  +5.14%  "

$ws.Range("D42").Value = "'1.56"
$ws.Range("E42").Value = "  +13.08%  "

$ws.Range("D43").Value = "'0.684"
$ws.Range("E43").Value = "  +5.06%  "

$ws.Range("D44").Value = "'3.146.95"
$ws.Range("E44").Value = "  +1.58%  "

$ws.Range("D45").Value = "'36.84"

$ws.Range("D46").Value = "'0.999"

$ws.Range("E49").Value = "  +8.19%  "

$ws.Range("D50").Value = "'20.72"
$ws.Range("E50").Value = "  +1.43%  "

$ws.Range("E51").Value = "  +3.58%  "

# --------------------------------------------------------------------------
# Rows 47 and 48 swap ranking positions in this snapshot: VeChain (previously
# row 47) drops to row 48, and Maker (previously row 48) rises to row 47.
# Update every column (Coin, Link, Price, Volume) to reflect the new order.
# --------------------------------------------------------------------------
$ws.Range("B47").Value = "Maker"
$ws.Range("C47").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D47").Value = "'2.293.04"
$ws.Range("E47").Value = "  +1.03%  "

$ws.Range("B48").Value = "VeChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D48").Value = "'0.0262"
$ws.Range("E48").Value = "  +3.71%  "

# --------------------------------------------------------------------------
# The leading apostrophe above stamps each touched Price cell with a
# "quote prefix" text style. Reset the whole Price column back to the
# workbook's plain default style so formatting matches the rest of the
# (unstyled) data cells - untouched cells in the range are unaffected.
# --------------------------------------------------------------------------
$ws.Range("D2:D51").Style = "Normal"
